# Commit: "removed test chapter. finalized matrix. added content to comparrison."
#
# Semantic changes being applied to Sheet1 ("Matrix"):
#  - Consolidate the subjective-rating vocabulary: the separate "gering" /
#    "mittel" ratings collapse into a single "niedrig" rating, and the
#    "keine Produkte am Markt" placeholder in the objective table is
#    replaced by a real measured value (720).
#  - Drop the gray highlight fill that was used as a scratch/"test" marker
#    on a handful of cells (now finalized, no longer needs to stand out).
#  - Row 6 no longer needs the taller wrapped height now that B6 holds a
#    short number instead of long text; row 9 reverts to the default row
#    height.
#  - Reset the view: scroll back to the top and select D5 instead of G14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Content edits
# ---------------------------------------------------------------------

# Objective table: "keine Produkte am Markt" -> measured value 720
$ws.Range("B6").Value = 720

# Subjective table: consolidate "gering"/"mittel" -> "niedrig", and the
# lone "sehr hoch" outlier -> "hoch"
$ws.Range("C10").Value = "niedrig"
$ws.Range("C11").Value = "niedrig"
$ws.Range("C12").Value = "niedrig"
$ws.Range("D12").Value = "niedrig"
$ws.Range("C14").Value = "hoch"
$ws.Range("D14").Value = "niedrig"

# ---------------------------------------------------------------------
# 2. Formatting: remove the leftover gray highlight fill (finalized, no
#    longer a "work in progress" marker)
# ---------------------------------------------------------------------

$highlighted = @("B3", "B6", "C10", "C11", "C12", "C13", "C14")
foreach ($addr in $highlighted) {
    $ws.Range($addr).Interior.ColorIndex = -4142
    $ws.Range($addr).Interior.Pattern = -4142
}

# ---------------------------------------------------------------------
# 3. Row heights
# ---------------------------------------------------------------------

# Row 6 shrinks now that B6 is a short number instead of wrapped text
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(6).RowHeight = 45

# Row 9 goes back to the default (non-custom) row height
$ws.Rows.Item(9).AutoFit()

# ---------------------------------------------------------------------
# 4. View: scroll back to the top and select D5 (was G14, scrolled to A7)
# ---------------------------------------------------------------------

$ws.Range("D5").Select()
